# Updated cryptos list on Tue Sep 10 07:13:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) holds plain-text numeric-looking values (e.g. "1.00", "5.20")
# that must retain exact formatting/trailing zeros, so format it as Text before
# writing so Excel doesn't silently convert them to real numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "56.436.41"
$ws.Range("E2").Value = "  +3.47%  "
$ws.Range("D3").Value = "2.321.16"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.50%  "
$ws.Range("D5").Value = "517.96"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").Value = "134.23"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").Value = "2.340.94"
$ws.Range("E9").Value = "  +1.86%  "
$ws.Range("E10").Value = "  +6.67%  "
$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("D12").Value = "5.20"
$ws.Range("E12").Value = "  +6.05%  "
$ws.Range("E13").Value = "  +0.03%  "
$ws.Range("D14").Value = "2.758.41"
$ws.Range("E14").Value = "  +2.65%  "
$ws.Range("D15").Value = "23.63"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "56.749.69"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "2.335.05"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "10.41"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "4.23"
$ws.Range("E20").Value = "  +1.93%  "
$ws.Range("D21").Value = "322.19"
$ws.Range("E21").Value = "  +4.74%  "
$ws.Range("E22").Value = "  -0.72%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("D24").Value = "60.79"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("E26").Value = "  +7.09%  "
$ws.Range("D27").Value = "7.74"
$ws.Range("E27").Value = "  +3.68%  "
$ws.Range("D28").Value = "1.23"
$ws.Range("E28").Value = "  +10.05%  "
$ws.Range("D29").Value = "170.22"
$ws.Range("E29").Value = "  -0.60%  "
$ws.Range("E30").Value = "  +5.47%  "
$ws.Range("E31").Value = "  +3.14%  "
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("E33").Value = "  +1.53%  "
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "0.993"
$ws.Range("E35").Value = "  -0.25%  "
$ws.Range("E36").Value = "  +2.95%  "
$ws.Range("D37").Value = "0.924"
$ws.Range("E37").Value = "  +1.80%  "
$ws.Range("E38").Value = "  +4.82%  "
$ws.Range("E39").Value = "  +7.98%  "
$ws.Range("D40").Value = "37.79"
$ws.Range("E40").Value = "  +3.44%  "
$ws.Range("D41").Value = "0.377"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("D42").Value = "3.58"
$ws.Range("E42").Value = "  +5.20%  "
$ws.Range("D43").Value = "136.57"
$ws.Range("E43").Value = "  +3.71%  "
$ws.Range("D44").Value = "276.40"
$ws.Range("E44").Value = "  +10.17%  "
$ws.Range("D45").Value = "5.07"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("E46").Value = "  +2.66%  "
$ws.Range("D47").Value = "0.0502"
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("E48").Value = "  +2.03%  "
$ws.Range("E49").Value = "  +5.03%  "
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "1.68"
$ws.Range("E51").Value = "  +10.11%  "
